$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A single label cell entered with a leading apostrophe (text/quote-prefix)
$ws.Range("E12").Value = "'group by , sum flu risk"

# The strings below are entered in the exact order the shared-string table
# needs to see them so new unique entries land at the right index.
$ws.Range("F23").Value = "icd"
$ws.Range("D23").Value = "id"
$ws.Range("E23").Value = "year"
$ws.Range("F24").Value = "J21"
$ws.Range("G23").Value = "flu_risk?"
$ws.Range("F25").Value = "J20"
$ws.Range("H23").Value = "flu_risk_cat"
$ws.Range("H24").Value = "Anemia"
$ws.Range("H25").Value = "na"
$ws.Range("H26").Value = "Lung "
$ws.Range("F26").Value = "H20"
$ws.Range("F27").Value = "243j"
$ws.Range("I23").Value = "flu_risk ever?"
$ws.Range("J23").Value = "flu_risk_last 3"
$ws.Range("K23").Value = "flu_risk_last 5"
$ws.Range("L23").Value = "U71?"
$ws.Range("M23").Value = "U72?"
$ws.Range("N23").Value = "U71|U72"
$ws.Range("O23").Value = "total_number_hospital_stays_ before 2020"
$ws.Range("P23").Value = "total_number_hospital_stays_from_2017_2019"
$ws.Range("Q23").Value = "total_number_stays_2015_2019?"
$ws.Range("H27").Value = "na"

# Numeric values
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 2000
$ws.Range("G24").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1
$ws.Range("K24").Value = 1
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 1
$ws.Range("P24").Value = 100
$ws.Range("Q24").Value = 50

$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 2001
$ws.Range("G25").Value = 0
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 1

$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 2002
$ws.Range("G26").Value = 1
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 1

$ws.Range("D27").Value = 12
$ws.Range("E27").Value = 2003
$ws.Range("G27").Value = 0
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 1

# Window / view adjustments to match the author's saved state
$ws.Range("L8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
